$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update values ---
$ws.Cells.Item(2, 1).Value = 'Colombia'
$ws.Cells.Item(2, 2).NumberFormat = "@"
$ws.Cells.Item(2, 2).Value = '1'
$ws.Cells.Item(2, 3).Value = 'Drugs (Pharmaceutical)'
$ws.Cells.Item(2, 7).Value = -168
$ws.Cells.Item(2, 8).Value = -168
$ws.Cells.Item(2, 9).Value = -121
$ws.Cells.Item(2, 10).Value = -121
$ws.Cells.Item(2, 11).Value = -3.82
$ws.Cells.Item(2, 12).Value = -127.3333333333333
$ws.Cells.Item(2, 13).Value = 0
$ws.Cells.Item(2, 14).Value = 0
$ws.Cells.Item(2, 15).Value = -0
$ws.Cells.Item(2, 16).Value = 0
$ws.Cells.Item(2, 17).Value = 0
$ws.Cells.Item(2, 18).Value = -0
$ws.Cells.Item(2, 19).Value = 0
$ws.Cells.Item(2, 21).Value = 0.554
$ws.Cells.Item(2, 22).Value = 0.0557344064386318
$ws.Cells.Item(2, 23).Value = -0.5701492537313433
$ws.Cells.Item(2, 24).Value = 0.07471565761744819
$ws.Cells.Item(2, 25).Value = -0.6448649113487914
$ws.Cells.Item(2, 26).Value = 0.008973975471133711
$ws.Cells.Item(2, 27).Value = -1.085851032007179
$ws.Cells.Item(2, 28).Value = 0.07037233853248395
$ws.Cells.Item(2, 29).Value = -1.156223370539663
$ws.Cells.Item(2, 30).Value = 1.12
$ws.Cells.Item(2, 31).Value = 0
$ws.Cells.Item(2, 32).Value = 1.12
$ws.Cells.Item(2, 33).Value = 0.5660000000000001
$ws.Cells.Item(2, 34).Value = 0.1012658227848102
$ws.Cells.Item(2, 35).Value = 0.2196078431372549
$ws.Cells.Item(2, 36).Value = 0.05387397677517609
$ws.Cells.Item(2, 37).Value = 0.1245050593928729
$ws.Cells.Item(2, 38).Value = 0.034
$ws.Cells.Item(2, 39).Value = 0.034
$ws.Cells.Item(2, 40).Value = -0.3522012578616353
$ws.Cells.Item(2, 41).Value = -106.7647058823529
$ws.Cells.Item(2, 42).Value = -0.1779874213836478
$ws.Cells.Item(2, 43).Value = -106.7647058823529

# --- Row 3: update values (now Blueberries Medical Corp.) ---
$ws.Cells.Item(3, 1).Value = 'Colombia'
$ws.Cells.Item(3, 2).Value = 'Blueberries Medical Corp. (CNSX:BBM)'
$ws.Cells.Item(3, 3).Value = 'Drugs (Pharmaceutical)'
$ws.Cells.Item(3, 7).Value = -168
$ws.Cells.Item(3, 8).Value = -168
$ws.Cells.Item(3, 9).Value = -121
$ws.Cells.Item(3, 10).Value = -121
$ws.Cells.Item(3, 11).Value = -3.82
$ws.Cells.Item(3, 12).Value = -127.3333333333333
$ws.Cells.Item(3, 13).Value = -0
$ws.Cells.Item(3, 14).Value = -0
$ws.Cells.Item(3, 15).Value = 0
$ws.Cells.Item(3, 16).Value = -0
$ws.Cells.Item(3, 17).Value = -0
$ws.Cells.Item(3, 18).Value = 0
$ws.Cells.Item(3, 19).Value = 0
$ws.Cells.Item(3, 21).Value = 0.554
$ws.Cells.Item(3, 22).Value = 0.0557344064386318
$ws.Cells.Item(3, 23).Value = -0.5701492537313433
$ws.Cells.Item(3, 24).Value = 0.07471565761744819
$ws.Cells.Item(3, 25).Value = -0.6448649113487914
$ws.Cells.Item(3, 26).Value = 0.008973975471133711
$ws.Cells.Item(3, 27).Value = -1.085851032007179
$ws.Cells.Item(3, 28).Value = 0.07037233853248395
$ws.Cells.Item(3, 29).Value = -1.156223370539663
$ws.Cells.Item(3, 30).Value = 1.12
$ws.Cells.Item(3, 31).Value = 0
$ws.Cells.Item(3, 32).Value = 1.12
$ws.Cells.Item(3, 33).Value = 0.5660000000000001
$ws.Cells.Item(3, 34).Value = 0.1012658227848102
$ws.Cells.Item(3, 35).Value = 0.2196078431372549
$ws.Cells.Item(3, 36).Value = 0.05387397677517609
$ws.Cells.Item(3, 37).Value = 0.1245050593928729
$ws.Cells.Item(3, 38).Value = 0.034
$ws.Cells.Item(3, 39).Value = 0.034
$ws.Cells.Item(3, 40).Value = -0.3522012578616353
$ws.Cells.Item(3, 41).Value = -106.7647058823529
$ws.Cells.Item(3, 42).Value = -0.1779874213836478
$ws.Cells.Item(3, 43).Value = -106.7647058823529

# --- Remove old row 4 (data merged into row 3; table now ends at row 3) ---
$ws.Rows(4).Delete()
